# Check.xlsx update:
#  - Fill in the "Отметка о прохождении" (column D) on sheet 1 with Pass/Fail marks
#    for every existing check row (2-27).
#  - Append a new check group "Нестиабильное/отсутствие интернет-соединения" as
#    rows 28-29 (merged A28:A29), with its own Pass marks in column D.
#  - Extend the list data validation on column C down to the new rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1) Give column D (rows 2-27) the same plain bordered look as columns B/C ---
# (it previously used the "empty" style with no side borders filled in)
$ws.Range("C2").Copy()
$ws.Range("D2:D27").PasteSpecial(-4122)

# --- 2) Write the Pass/Fail marks for the existing rows ---
$ws.Cells.Item(2, 4).Value = "Pass"
$ws.Cells.Item(3, 4).Value = "Pass"
$ws.Cells.Item(4, 4).Value = "Pass"
$ws.Cells.Item(5, 4).Value = "Fail"
$ws.Cells.Item(6, 4).Value = "Pass"
$ws.Cells.Item(7, 4).Value = "Pass"
$ws.Cells.Item(8, 4).Value = "Pass"
$ws.Cells.Item(9, 4).Value = "Pass"
$ws.Cells.Item(10, 4).Value = "Pass"
$ws.Cells.Item(11, 4).Value = "Pass"
$ws.Cells.Item(12, 4).Value = "Pass"
$ws.Cells.Item(13, 4).Value = "Pass"
$ws.Cells.Item(14, 4).Value = "Pass"
$ws.Cells.Item(15, 4).Value = "Fail"
$ws.Cells.Item(16, 4).Value = "Pass"
$ws.Cells.Item(17, 4).Value = "Pass"
$ws.Cells.Item(18, 4).Value = "Pass"
$ws.Cells.Item(19, 4).Value = "Pass"
$ws.Cells.Item(20, 4).Value = "Pass"
$ws.Cells.Item(21, 4).Value = "Pass"
$ws.Cells.Item(22, 4).Value = "Pass"
$ws.Cells.Item(23, 4).Value = "Pass"
$ws.Cells.Item(24, 4).Value = "Pass"
$ws.Cells.Item(25, 4).Value = "Pass"
$ws.Cells.Item(26, 4).Value = "Pass"
$ws.Cells.Item(27, 4).Value = "Fail"

# --- 3) Append the new "Нестиабильное/отсутствие интернет-соединения" group (rows 28-29) ---
# Merge A28:A29 first, then copy the formatting of another two-row merged group
# (rows 20-21) onto it so the borders/fills/fonts match the rest of the table
# (merging after the format paste would clobber the per-row font), then
# overwrite the text/values.
$ws.Range("A28:A29").Merge()

$ws.Range("A20:D21").Copy()
$ws.Range("A28:D29").PasteSpecial(-4122)

$ws.Cells.Item(28, 1).Value = "Нестиабильное/отсутствие интернет-соединения"
$ws.Cells.Item(28, 2).Value = "Авторизация при плохом интернет-соединении"
$ws.Cells.Item(28, 3).Value = "High"
$ws.Cells.Item(28, 4).Value = "Pass"

$ws.Cells.Item(29, 2).Value = "Авторизация при полном отсутствии интернет-соединения"
$ws.Cells.Item(29, 3).Value = "High"
$ws.Cells.Item(29, 4).Value = "Pass"

# Long group label wraps inside its cell
$ws.Range("A28").WrapText = $true

# --- 4) Extend the column C list validation down to the new rows ---
$ws.Range("C2:C29").Validation.Delete()
$ws.Range("C2:C29").Validation.Add(3, 1, 1, "High,Medium,Low")
